$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Append two new API call log rows to the "API Calls" sheet
# ---------------------------------------------------------------------
$apiCalls = $wb.Worksheets.Item("API Calls")

$apiCalls.Range("A14").Value = "2025-12-15 14:19:40"
$apiCalls.Range("B14").Value = "Remote/SSH Cline"
$apiCalls.Range("C14").Value = "Smoothie_Bar_IG"
$apiCalls.Range("D14").Value = "right now i have a generate_report script that takes a standard word doc and makes it into an HTML f"
$apiCalls.Range("E14").Value = "claude-sonnet-3-5"
$apiCalls.Range("F14").Value = 1506903
$apiCalls.Range("G14").Value = 62835
$apiCalls.Range("H14").Value = 13.435237
$apiCalls.Range("I14").Value = 24.016044

$apiCalls.Range("A15").Value = "2025-12-15 19:34:46"
$apiCalls.Range("B15").Value = "Remote/SSH Cline"
$apiCalls.Range("C15").Value = "Smoothie_Bar_IG"
$apiCalls.Range("D15").Value = "i am making a script that transforms this word doc as a template to html for presentation.   i want "
$apiCalls.Range("E15").Value = "claude-sonnet-3-5"
$apiCalls.Range("F15").Value = 1917669
$apiCalls.Range("G15").Value = 46976
$apiCalls.Range("H15").Value = 16.337973
$apiCalls.Range("I15").Value = 40.354017

# ---------------------------------------------------------------------
# 2. Add a new "Summary8" sheet after "Summary7", matching the layout /
#    formatting of the existing Summary sheets
# ---------------------------------------------------------------------
$srcSheet = $wb.Worksheets.Item("Summary7")
$afterSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $afterSheet)
$newSheet.Name = "Summary8"

# Copy formatting (styles, fills, fonts) from the previous summary sheet
$srcSheet.Range("A1:B10").Copy() | Out-Null
$newSheet.Range("A1").PasteSpecial(-4122) | Out-Null

# Match column widths
$newSheet.Columns.Item(1).ColumnWidth = $srcSheet.Columns.Item(1).ColumnWidth
$newSheet.Columns.Item(2).ColumnWidth = $srcSheet.Columns.Item(2).ColumnWidth

$newSheet.Range("A1").Value = "Metric"
$newSheet.Range("B1").Value = "Value"

$newSheet.Range("A2").Value = "Total API Calls"
$newSheet.Range("B2").Value = 2

$newSheet.Range("A3").Value = "- Cline Coding Calls"
$newSheet.Range("B3").Value = 2

$newSheet.Range("A4").Value = "- Automated Analysis Calls"
$newSheet.Range("B4").Value = 0

$newSheet.Range("A5").Value = "Average Cost per Call (All)"
$newSheet.Range("A6").Value = "- Avg Coding Call Cost"
$newSheet.Range("A7").Value = "- Avg Analysis Call Cost"
$newSheet.Range("A8").Value = "Total Cost"
$newSheet.Range("A9").Value = "- Coding Cost"
$newSheet.Range("A10").Value = "- Analysis Cost"

# These look like currency amounts ("$12.34"); Excel would normally auto-
# convert such strings into formatted numbers when assigned directly, so
# force the cells to Text format first, assign the literal string, then
# re-apply the original (non-text) cell formatting copied from the
# template sheet so the visual style still matches the other summary
# sheets while the stored value remains the literal text.
$moneyCells = @("B5", "B6", "B7", "B8", "B9", "B10")
$moneyValues = @("$20.177009", "$14.886605", "$0.000000", "$40.354017", "$29.773210", "$0.000000")
for ($i = 0; $i -lt $moneyCells.Length; $i++) {
    $cellAddr = $moneyCells[$i]
    $newSheet.Range($cellAddr).NumberFormat = "@"
    $newSheet.Range($cellAddr).Value = $moneyValues[$i]
    $srcSheet.Range($cellAddr).Copy() | Out-Null
    $newSheet.Range($cellAddr).PasteSpecial(-4122) | Out-Null
}

# Restore the originally active sheet/selection so the workbook view
# isn't left pointing at the newly created sheet.
$apiCalls.Activate()
$apiCalls.Range("A1").Select() | Out-Null
